$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("A2:C92")
$key = $ws.Range("A1")
$rng.Sort($key, 1)
